$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("wenden", "none", "none")
    ,@("danken", "flower/flower014.jpg", "flower")
    ,@("schicken", "flower/flower003.jpg", "flower")
    ,@("rufen", "none", "none")
    ,@("bieten", "flower/flower008.jpg", "flower")
    ,@("tropfen", "flower/flower029.jpg", "flower")
    ,@("enden", "none", "none")
    ,@("bilden", "flower/flower019.jpg", "flower")
    ,@("planen", "dog/dog021.jpg", "dog")
    ,@("opfern", "none", "none")
    ,@("kleben", "flower/flower017.jpg", "flower")
    ,@("geben", "flower/flower007.jpg", "flower")
    ,@("schalten", "none", "none")
    ,@("fühlen", "flower/flower024.jpg", "flower")
    ,@("dringen", "dog/dog006.jpg", "dog")
    ,@("klagen", "none", "none")
    ,@("lassen", "dog/dog001.jpg", "dog")
    ,@("seufzen", "flower/flower010.jpg", "flower")
    ,@("sparen", "none", "none")
    ,@("sammeln", "flower/flower011.jpg", "flower")
    ,@("posten", "dog/dog030.jpg", "dog")
    ,@("hören", "none", "none")
    ,@("schultern", "flower/flower033.jpg", "flower")
    ,@("grüßen", "dog/dog024.jpg", "dog")
    ,@("drohen", "none", "none")
    ,@("bitten", "flower/flower001.jpg", "flower")
    ,@("kehren", "dog/dog023.jpg", "dog")
    ,@("orten", "none", "none")
    ,@("zeugen", "dog/dog026.jpg", "dog")
    ,@("zahlen", "dog/dog007.jpg", "dog")
    ,@("weigern", "none", "none")
    ,@("beißen", "dog/dog022.jpg", "dog")
    ,@("altern", "dog/dog010.jpg", "dog")
    ,@("ächzen", "none", "none")
    ,@("öffnen", "dog/dog011.jpg", "dog")
    ,@("drücken", "flower/flower025.jpg", "flower")
    ,@("kosten", "none", "none")
    ,@("runden", "flower/flower021.jpg", "flower")
    ,@("gleichen", "dog/dog012.jpg", "dog")
    ,@("dauern", "none", "none")
    ,@("schwimmen", "flower/flower032.jpg", "flower")
    ,@("pflegen", "dog/dog028.jpg", "dog")
    ,@("stören", "none", "none")
    ,@("ändern", "dog/dog000.jpg", "dog")
    ,@("fordern", "dog/dog027.jpg", "dog")
    ,@("bremsen", "none", "none")
    ,@("münzen", "flower/flower015.jpg", "flower")
    ,@("wüten", "dog/dog025.jpg", "dog")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}